# final updates for sending
#
# - Normalize the "Genero" column (C) values from "M"/"F" to lowercase "m"/"f"
# - Update the selected cell on the sheet to C30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Formula
    if ($current -eq "M") {
        $cell.Value = "m"
    } elseif ($current -eq "F") {
        $cell.Value = "f"
    }
}

$ws.Range("C30").Select()
